# Residencial plans update (test(web)/qa/salesforce/residencial - agregar planes de residencial)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plans")

# Row 5 (Plan 1): Tipo de Servicio becomes Sin_TotalPlay_TV, Amazon Prime cleared
$ws.Range("C5").Value = "Sin_TotalPlay_TV"
$ws.Range("E5").Value = ""

# Row 6 (Plan 2): Tipo de Servicio becomes Sin_TotalPlay_TV, Megas 50 -> 100, Netflix cleared
$ws.Range("C6").Value = "Sin_TotalPlay_TV"
$ws.Range("D6").Value = 100
$ws.Range("F6").Value = ""

# Row 7 (Plan 3): Tipo de Servicio becomes Sin_TotalPlay_TV, Megas 50 -> 200, Amazon Prime cleared
$ws.Range("C7").Value = "Sin_TotalPlay_TV"
$ws.Range("D7").Value = 200
$ws.Range("E7").Value = ""

# Row 8 (Plan 4): Tipo de Servicio becomes Sin_TotalPlay_TV, Megas 50 -> 500, Netflix cleared
$ws.Range("C8").Value = "Sin_TotalPlay_TV"
$ws.Range("D8").Value = 500
$ws.Range("F8").Value = ""

# Update active selection to reflect the last-edited cell in the Plans sheet
$ws.Range("D8").Select()
